$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates ---
$ws.Range("A8").Value = "Volume 32   Number  46"
$ws.Range("C9").Value = "Report Covering the Week  11/10/2025  Through  11/16/2025"

# --- Fix cell styles for cells that change between numeric and placeholder-text ---
$ws.Range("I29").Copy()
$ws.Range("C20").PasteSpecial(-4122)
$ws.Range("I29").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("K29").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("I29").Copy()
$ws.Range("D33").PasteSpecial(-4122)
$ws.Range("K29").Copy()
$ws.Range("E33").PasteSpecial(-4122)
$ws.Range("I29").Copy()
$ws.Range("G33").PasteSpecial(-4122)
$ws.Range("K29").Copy()
$ws.Range("H33").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- C22: numeric 1 -> placeholder text "0" (style 13), using a forced-text formula then freezing to a value ---
$ws.Range("C29").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("C22").Formula = "=""0"""
$ws.Range("C22").Copy()
$ws.Range("C22").PasteSpecial(-4163)
$excel.CutCopyMode = $false

# --- Numeric value updates ---
# Row 15
$ws.Range("D15").Value = 3
$ws.Range("G15").Value = 5
$ws.Range("H15").Value = -60
$ws.Range("I15").Value = 20
$ws.Range("J15").Value = 24
$ws.Range("K15").Value = -16.666666666666
$ws.Range("L15").Value = 33.333333333333
$ws.Range("M15").Value = 11.111111111111
$ws.Range("N15").Value = -42.857142857142

# Row 16
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = -75
$ws.Range("F16").Value = 7
$ws.Range("G16").Value = 13
$ws.Range("H16").Value = -46.153846153846
$ws.Range("I16").Value = 92
$ws.Range("J16").Value = 169
$ws.Range("K16").Value = -45.562130177514
$ws.Range("L16").Value = -33.812949640287
$ws.Range("M16").Value = -57.009345794392
$ws.Range("N16").Value = -90.075512405609

# Row 17
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 7
$ws.Range("E17").Value = -28.571428571428
$ws.Range("F17").Value = 15
$ws.Range("G17").Value = 24
$ws.Range("H17").Value = -37.5
$ws.Range("I17").Value = 265
$ws.Range("J17").Value = 292
$ws.Range("K17").Value = -9.246575342465
$ws.Range("L17").Value = 4.743083003952
$ws.Range("M17").Value = 54.970760233918
$ws.Range("N17").Value = -43.1330472103

# Row 18
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 7
$ws.Range("G18").Value = 9
$ws.Range("H18").Value = -22.222222222222
$ws.Range("I18").Value = 92
$ws.Range("J18").Value = 149
$ws.Range("K18").Value = -38.255033557047
$ws.Range("L18").Value = -21.367521367521
$ws.Range("M18").Value = -58.928571428571
$ws.Range("N18").Value = -91.033138401559

# Row 19
$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 4
$ws.Range("E19").Value = 100
$ws.Range("F19").Value = 29
$ws.Range("G19").Value = 36
$ws.Range("H19").Value = -19.444444444444
$ws.Range("I19").Value = 310
$ws.Range("J19").Value = 427
$ws.Range("K19").Value = -27.400468384074
$ws.Range("L19").Value = -38.12375249501
$ws.Range("M19").Value = -5.487804878048
$ws.Range("N19").Value = -24.205378973105

# Row 20
$ws.Range("C20").Value = 4
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = 100
$ws.Range("F20").Value = 9
$ws.Range("G20").Value = 10
$ws.Range("H20").Value = -10
$ws.Range("I20").Value = 93
$ws.Range("J20").Value = 135
$ws.Range("K20").Value = -31.111111111111
$ws.Range("L20").Value = -41.875
$ws.Range("M20").Value = -14.678899082568
$ws.Range("N20").Value = -89.211136890951

# Row 21
$ws.Range("C21").Value = 20
$ws.Range("D21").Value = 22
$ws.Range("E21").Value = -9.090909090909
$ws.Range("F21").Value = 69
$ws.Range("G21").Value = 97
$ws.Range("H21").Value = -28.865979381443
$ws.Range("I21").Value = 873
$ws.Range("J21").Value = 1199
$ws.Range("K21").Value = -27.18932443703
$ws.Range("L21").Value = -26.515151515151
$ws.Range("M21").Value = -18.411214953271
$ws.Range("N21").Value = -76.670229823623

# Row 22
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = -100
$ws.Range("F22").Value = 2
$ws.Range("G22").Value = 1
$ws.Range("H22").Value = 100
$ws.Range("I22").Value = 31
$ws.Range("J22").Value = 22
$ws.Range("K22").Value = 40.90909090909
$ws.Range("L22").Value = -3.125
$ws.Range("M22").Value = 19.230769230769

# Row 24
$ws.Range("C24").Value = 15
$ws.Range("D24").Value = 16
$ws.Range("E24").Value = -6.25
$ws.Range("F24").Value = 63
$ws.Range("H24").Value = -22.222222222222
$ws.Range("I24").Value = 638
$ws.Range("J24").Value = 850
$ws.Range("K24").Value = -24.941176470588
$ws.Range("L24").Value = -42.157751586582
$ws.Range("M24").Value = -8.987161198288

# Row 25
$ws.Range("C25").Value = 4
$ws.Range("D25").Value = 7
$ws.Range("E25").Value = -42.857142857142
$ws.Range("F25").Value = 16
$ws.Range("H25").Value = -58.974358974359
$ws.Range("I25").Value = 187
$ws.Range("J25").Value = 349
$ws.Range("K25").Value = -46.418338108882
$ws.Range("L25").Value = -70.127795527156

# Row 26
$ws.Range("C26").Value = 9
$ws.Range("D26").Value = 11
$ws.Range("E26").Value = -18.181818181818
$ws.Range("F26").Value = 52
$ws.Range("G26").Value = 52
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 480
$ws.Range("J26").Value = 500
$ws.Range("K26").Value = -4
$ws.Range("L26").Value = 2.345415778251
$ws.Range("M26").Value = -14.438502673796

# Row 27
$ws.Range("D27").Value = 3
$ws.Range("F27").Value = 3
$ws.Range("G27").Value = 5
$ws.Range("H27").Value = -40
$ws.Range("J27").Value = 30
$ws.Range("K27").Value = -6.666666666666
$ws.Range("L27").Value = 21.739130434782

# Row 28
$ws.Range("F28").Value = 5
$ws.Range("H28").Value = -28.571428571428
$ws.Range("J28").Value = 51
$ws.Range("K28").Value = -3.92156862745
$ws.Range("L28").Value = -34.666666666666

# Row 33
$ws.Range("D33").Value = 1
$ws.Range("E33").Value = -100
$ws.Range("G33").Value = 1
$ws.Range("H33").Value = -100
$ws.Range("J33").Value = 5
$ws.Range("K33").Value = -20

